$d = $word.ActiveDocument

$replacements = @(
    @("78÷4=", "91÷3="),
    @("58÷4=", "40÷4="),
    @("24÷8=", "87÷8="),
    @("35÷7=", "70÷8="),
    @("53÷3=", "94÷9="),
    @("62÷4=", "17÷4="),
    @("18÷3=", "93÷9="),
    @("19÷8=", "66÷8="),
    @("12÷7=", "28÷5="),
    @("30÷3=", "52÷2="),
    @("63÷7=", "51÷7="),
    @("34÷7=", "85÷4="),
    @("10÷8=", "61÷9="),
    @("27÷9=", "10÷6="),
    @("95÷4=", "86÷8="),
    @("71÷4=", "86÷7="),
    @("53÷7=", "19÷8="),
    @("69÷4=", "14÷4="),
    @("85÷2=", "87÷2="),
    @("53÷2=", "70÷4="),
    @("89÷3=", "80÷6="),
    @("96÷3=", "73÷4="),
    @("68÷2=", "90÷5="),
    @("66÷9=", "11÷6="),
    @("14÷7=", "10÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "All replacements applied"
